$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3111.8572
$ws.Range("I6").Value = 6694.3335
$ws.Range("J6").Value = 425
$ws.Range("K6").Value = 20083.0005
$ws.Range("L6").Value = 1275
$ws.Range("M6").Value = -19971.0005
$ws.Range("N6").Value = -1499
$ws.Range("H15").Value = 1608.3448
$ws.Range("I15").Value = 1608.3448
$ws.Range("K15").Value = 4825.0344
$ws.Range("M15").Value = -4656.0344
$ws.Range("H33").Value = 476.11765
$ws.Range("I33").Value = 243.1
$ws.Range("K33").Value = 243.1
$ws.Range("M33").Value = -14.09999999999999
$ws.Range("H40").Value = 1251
$ws.Range("I40").Value = 1190
$ws.Range("J40").Value = 1266.25
$ws.Range("K40").Value = 1190
$ws.Range("L40").Value = 1266.25
$ws.Range("M40").Value = -1015
$ws.Range("N40").Value = -1616.25
$ws.Range("H69").Value = 2401.3333
$ws.Range("I69").Value = 2478.25
$ws.Range("J69").Value = 2247.5
$ws.Range("K69").Value = 7434.75
$ws.Range("L69").Value = 6742.5
$ws.Range("M69").Value = -6560.75
$ws.Range("N69").Value = -8490.5
$ws.Range("H72").Value = 2401.3333
$ws.Range("I72").Value = 2478.25
$ws.Range("J72").Value = 2247.5
$ws.Range("K72").Value = 22304.25
$ws.Range("L72").Value = 20227.5
$ws.Range("M72").Value = -17936.25
$ws.Range("N72").Value = -28963.5
$ws.Range("H112").Value = 1798.5333
$ws.Range("J112").Value = 2092.9
$ws.Range("L112").Value = 6278.700000000001
$ws.Range("N112").Value = -8494.700000000001
$ws.Range("H116").Value = 6817.6665
$ws.Range("I116").Value = 5770
$ws.Range("K116").Value = 5770
$ws.Range("M116").Value = -2328
$ws.Range("H127").Value = 1895.4546
$ws.Range("I127").Value = 1385
$ws.Range("K127").Value = 4155
$ws.Range("M127").Value = 805
$ws.Range("H132").Value = 2562.5454
$ws.Range("I132").Value = 2562.5454
$ws.Range("K132").Value = 7687.6362
$ws.Range("M132").Value = -5157.6362
$ws.Range("H138").Value = 942.8125
$ws.Range("J138").Value = 2000
$ws.Range("L138").Value = 6000
$ws.Range("N138").Value = -16280

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2262.4285
$ws.Range("I2").Value = 2262.4285
$ws.Range("K2").Value = 2262.4285
$ws.Range("M2").Value = -2149.4285
$ws.Range("H74").Value = 1006.6667
$ws.Range("I74").Value = 1006.6667
$ws.Range("K74").Value = 1006.6667
$ws.Range("M74").Value = -132.6667
$ws.Range("H77").Value = 1006.6667
$ws.Range("I77").Value = 1006.6667
$ws.Range("K77").Value = 5033.3335
$ws.Range("M77").Value = -665.3334999999997
$ws.Range("H116").Value = 2262.4285
$ws.Range("I116").Value = 2262.4285
$ws.Range("K116").Value = 2262.4285
$ws.Range("M116").Value = 31.57150000000001
$ws.Range("H132").Value = 1738.7333
$ws.Range("I132").Value = 1738.7333
$ws.Range("K132").Value = 5216.199900000001
$ws.Range("M132").Value = -2686.199900000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2262.4285
$ws.Range("I3").Value = 2262.4285
$ws.Range("K3").Value = 2262.4285
$ws.Range("M3").Value = -2148.4285
$ws.Range("H20").Value = 5291.1665
$ws.Range("I20").Value = 4349.4
$ws.Range("K20").Value = 4349.4
$ws.Range("M20").Value = -4102.4
$ws.Range("H86").Value = 2635.9092
$ws.Range("I86").Value = 2875.75
$ws.Range("J86").Value = 1996.3334
$ws.Range("K86").Value = 2875.75
$ws.Range("L86").Value = 1996.3334
$ws.Range("M86").Value = -1752.75
$ws.Range("N86").Value = -4242.3334
$ws.Range("H89").Value = 2635.9092
$ws.Range("I89").Value = 2875.75
$ws.Range("J89").Value = 1996.3334
$ws.Range("K89").Value = 14378.75
$ws.Range("L89").Value = 9981.666999999999
$ws.Range("M89").Value = -8762.75
$ws.Range("N89").Value = -21213.667
$ws.Range("H94").Value = 1531.6522
$ws.Range("I94").Value = 1773.2778
$ws.Range("J94").Value = 661.8
$ws.Range("K94").Value = 1773.2778
$ws.Range("L94").Value = 661.8
$ws.Range("M94").Value = -1322.2778
$ws.Range("N94").Value = -1563.8
$ws.Range("H107").Value = 873.2
$ws.Range("I107").Value = 873.2
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 873.2
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1046.8
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 6199.9565
$ws.Range("I134").Value = 7075.2144
$ws.Range("J134").Value = 4838.4443
$ws.Range("K134").Value = 21225.6432
$ws.Range("L134").Value = 14515.3329
$ws.Range("M134").Value = -18690.6432
$ws.Range("N134").Value = -19585.3329

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H132").Value = 2999.3333
$ws.Range("I132").Value = 2749.5
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 8248.5
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -5718.5
$ws.Range("N132").Value = -15557
$ws.Range("H134").Value = 3750.0667
$ws.Range("I134").Value = 3039.818
$ws.Range("K134").Value = 9119.454000000002
$ws.Range("M134").Value = -6584.454000000002

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6250.3335
$ws.Range("I70").Value = 4499.6665
$ws.Range("K70").Value = 4499.6665
$ws.Range("M70").Value = -4229.6665
$ws.Range("H73").Value = 6250.3335
$ws.Range("I73").Value = 4499.6665
$ws.Range("K73").Value = 4499.6665
$ws.Range("M73").Value = -3563.6665
$ws.Range("H80").Value = 2669.7273
$ws.Range("I80").Value = 994.5
$ws.Range("J80").Value = 3042
$ws.Range("K80").Value = 994.5
$ws.Range("L80").Value = 3042
$ws.Range("M80").Value = 3.5
$ws.Range("N80").Value = -5038
$ws.Range("H83").Value = 2669.7273
$ws.Range("I83").Value = 994.5
$ws.Range("J83").Value = 3042
$ws.Range("K83").Value = 4972.5
$ws.Range("L83").Value = 15210
$ws.Range("M83").Value = 19.5
$ws.Range("N83").Value = -25194
$ws.Range("H97").Value = 1108.9166
$ws.Range("I97").Value = 983
$ws.Range("J97").Value = 1234.8334
$ws.Range("K97").Value = 983
$ws.Range("L97").Value = 1234.8334
$ws.Range("M97").Value = -487
$ws.Range("N97").Value = -2226.8334

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2299
$ws.Range("I7").Value = 1748.75
$ws.Range("K7").Value = 1748.75
$ws.Range("M7").Value = -1636.75
$ws.Range("H40").Value = 6650.25
$ws.Range("I40").Value = 3300.5
$ws.Range("K40").Value = 3300.5
$ws.Range("M40").Value = -3164.5
$ws.Range("H61").Value = 9376.799999999999
$ws.Range("I61").Value = 9961.333000000001
$ws.Range("K61").Value = 9961.333000000001
$ws.Range("M61").Value = -9759.333000000001
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40450
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41560
$ws.Range("H113").Value = 9376.799999999999
$ws.Range("I113").Value = 9961.333000000001
$ws.Range("K113").Value = 9961.333000000001
$ws.Range("M113").Value = -7791.333000000001
$ws.Range("H126").Value = 2299
$ws.Range("I126").Value = 1748.75
$ws.Range("K126").Value = 5246.25
$ws.Range("M126").Value = -2776.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41248
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126240

Write-Output "Applied all profit-sheet updates"
